$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.748.18'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '2.493.06'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('E10').Value = '  +6.61%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.08%  '
$ws.Range('D15').Value = '2.882.67'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = '2.502.74'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.852'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '47.656.27'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('D21').Value = '0.0₃0942'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +13.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '247.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.139'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('E33').Value = '  +1.39%  '
$ws.Range('E34').Value = '  -2.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0792'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.112'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.25'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.99%  '
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').Value = '2.000.17'
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.81'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.84'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.24%  '
